$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values in row 1 for columns P (16) and Q (17)
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Match the formatting already used across the rest of row 1 (bold, bordered, centered)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats

for ($r = 2; $r -le 25; $r++) {
    # Flip values in columns I, K, M, O
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1

    # New columns P, Q
    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
